$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Higher Or Lower" column (D) values for each data row (rows 5-60)
# are no longer populated in this dataset - clear them out so the
# unused "-"/"+" shared strings are dropped and the remaining columns
# (E:O) keep their data as-is.
$ws.Range("D5:D60").ClearContents()
